$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 20,20
$arr[0,0] = "ECs"
$arr[0,1] = "Fn1"
$arr[0,2] = "Sdc2"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 66.47695399999999
$arr[0,7] = 199.430862
$arr[0,8] = 0.04311983106164722
$arr[0,9] = 0.04311983106164721
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 1.826566
$arr[0,13] = 5.479698
$arr[0,14] = 0.02795372904983374
$arr[0,15] = 0.02795372904983374
$arr[0,16] = 121.424543959964
$arr[0,17] = 1092.820895639676
$arr[0,18] = 0.001205360074171891
$arr[0,19] = 0.001205360074171891
$arr[1,0] = "ECs"
$arr[1,1] = "Fn1"
$arr[1,2] = "Sdc2"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 66.47695399999999
$arr[1,7] = 199.430862
$arr[1,8] = 0.04311983106164722
$arr[1,9] = 0.04311983106164721
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 44.29005966666667
$arr[1,13] = 132.870179
$arr[1,14] = 0.6778141756295529
$arr[1,15] = 0.6778141756295529
$arr[1,16] = 2944.268259118255
$arr[1,17] = 26498.4143320643
$arr[1,18] = 0.029227232744336
$arr[1,19] = 0.02922723274433599
$arr[2,0] = "ECs"
$arr[2,1] = "Fn1"
$arr[2,2] = "Sdc2"
$arr[2,3] = "MuSCs"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 66.47695399999999
$arr[2,7] = 199.430862
$arr[2,8] = 0.04311983106164722
$arr[2,9] = 0.04311983106164721
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 19.10886933333333
$arr[2,13] = 57.326608
$arr[2,14] = 0.2924417490485847
$arr[2,15] = 0.2924417490485847
$arr[2,16] = 1270.299427664011
$arr[2,17] = 11432.6948489761
$arr[2,18] = 0.0126100388143476
$arr[2,19] = 0.0126100388143476
$arr[3,0] = "ECs"
$arr[3,1] = "Fn1"
$arr[3,2] = "Sdc2"
$arr[3,3] = "Resolving-Mac"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 66.47695399999999
$arr[3,7] = 199.430862
$arr[3,8] = 0.04311983106164722
$arr[3,9] = 0.04311983106164721
$arr[3,10] = 2
$arr[3,11] = 0.6666666666666666
$arr[3,12] = 0.1169856666666667
$arr[3,13] = 0.350957
$arr[3,14] = 0.001790346272028586
$arr[3,15] = 0.001790346272028586
$arr[3,16] = 7.776850781659332
$arr[3,17] = 69.99165703493399
$arr[3,18] = 0.00007719942879172254
$arr[3,19] = 0.00007719942879172254
$arr[4,0] = "FAPs"
$arr[4,1] = "Fn1"
$arr[4,2] = "Sdc2"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 1361.379069
$arr[4,7] = 4084.137207
$arr[4,8] = 0.8830494168872806
$arr[4,9] = 0.8830494168872804
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 1.826566
$arr[4,13] = 5.479698
$arr[4,14] = 0.02795372904983374
$arr[4,15] = 0.02795372904983374
$arr[4,16] = 2486.648720547053
$arr[4,17] = 22379.83848492349
$arr[4,18] = 0.02468452413728072
$arr[4,19] = 0.02468452413728072
$arr[5,0] = "FAPs"
$arr[5,1] = "Fn1"
$arr[5,2] = "Sdc2"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 1361.379069
$arr[5,7] = 4084.137207
$arr[5,8] = 0.8830494168872806
$arr[5,9] = 0.8830494168872804
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 44.29005966666667
$arr[5,13] = 132.870179
$arr[5,14] = 0.6778141756295529
$arr[5,15] = 0.6778141756295529
$arr[5,16] = 60295.56019496112
$arr[5,17] = 542660.0417546501
$arr[5,18] = 0.5985434125476095
$arr[5,19] = 0.5985434125476095
$arr[6,0] = "FAPs"
$arr[6,1] = "Fn1"
$arr[6,2] = "Sdc2"
$arr[6,3] = "MuSCs"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 1361.379069
$arr[6,7] = 4084.137207
$arr[6,8] = 0.8830494168872806
$arr[6,9] = 0.8830494168872804
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 19.10886933333333
$arr[6,13] = 57.326608
$arr[6,14] = 0.2924417490485847
$arr[6,15] = 0.2924417490485847
$arr[6,16] = 26014.41474265598
$arr[6,17] = 234129.7326839039
$arr[6,18] = 0.2582405159708492
$arr[6,19] = 0.2582405159708491
$arr[7,0] = "FAPs"
$arr[7,1] = "Fn1"
$arr[7,2] = "Sdc2"
$arr[7,3] = "Resolving-Mac"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 1361.379069
$arr[7,7] = 4084.137207
$arr[7,8] = 0.8830494168872806
$arr[7,9] = 0.8830494168872804
$arr[7,10] = 2
$arr[7,11] = 0.6666666666666666
$arr[7,12] = 0.1169856666666667
$arr[7,13] = 0.350957
$arr[7,14] = 0.001790346272028586
$arr[7,15] = 0.001790346272028586
$arr[7,16] = 159.261837973011
$arr[7,17] = 1433.356541757099
$arr[7,18] = 0.00158096423154116
$arr[7,19] = 0.00158096423154116
$arr[8,0] = "Inflammatory-Mac"
$arr[8,1] = "Fn1"
$arr[8,2] = "Sdc2"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 44.831112
$arr[8,7] = 134.493336
$arr[8,8] = 0.02907940059566787
$arr[8,9] = 0.02907940059566786
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 1.826566
$arr[8,13] = 5.479698
$arr[8,14] = 0.02795372904983374
$arr[8,15] = 0.02795372904983374
$arr[8,16] = 81.88698492139198
$arr[8,17] = 736.982864292528
$arr[8,18] = 0.0008128776851828736
$arr[8,19] = 0.0008128776851828735
$arr[9,0] = "Inflammatory-Mac"
$arr[9,1] = "Fn1"
$arr[9,2] = "Sdc2"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 44.831112
$arr[9,7] = 134.493336
$arr[9,8] = 0.02907940059566787
$arr[9,9] = 0.02907940059566786
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 44.29005966666667
$arr[9,13] = 132.870179
$arr[9,14] = 0.6778141756295529
$arr[9,15] = 0.6778141756295529
$arr[9,16] = 1985.572625403016
$arr[9,17] = 17870.15362862714
$arr[9,18] = 0.01971042994255415
$arr[9,19] = 0.01971042994255414
$arr[10,0] = "Inflammatory-Mac"
$arr[10,1] = "Fn1"
$arr[10,2] = "Sdc2"
$arr[10,3] = "MuSCs"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 44.831112
$arr[10,7] = 134.493336
$arr[10,8] = 0.02907940059566787
$arr[10,9] = 0.02907940059566786
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 19.10886933333333
$arr[10,13] = 57.326608
$arr[10,14] = 0.2924417490485847
$arr[10,15] = 0.2924417490485847
$arr[10,16] = 856.671861276032
$arr[10,17] = 7710.046751484288
$arr[10,18] = 0.008504030771481569
$arr[10,19] = 0.008504030771481567
$arr[11,0] = "Inflammatory-Mac"
$arr[11,1] = "Fn1"
$arr[11,2] = "Sdc2"
$arr[11,3] = "Resolving-Mac"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 44.831112
$arr[11,7] = 134.493336
$arr[11,8] = 0.02907940059566787
$arr[11,9] = 0.02907940059566786
$arr[11,10] = 2
$arr[11,11] = 0.6666666666666666
$arr[11,12] = 0.1169856666666667
$arr[11,13] = 0.350957
$arr[11,14] = 0.001790346272028586
$arr[11,15] = 0.001790346272028586
$arr[11,16] = 5.244597524727999
$arr[11,17] = 47.201377722552
$arr[11,18] = 0.00005206219644927982
$arr[11,19] = 0.00005206219644927981
$arr[12,0] = "MuSCs"
$arr[12,1] = "Fn1"
$arr[12,2] = "Sdc2"
$arr[12,3] = "ECs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 52.83062100000001
$arr[12,7] = 158.491863
$arr[12,8] = 0.0342682285413064
$arr[12,9] = 0.03426822854130639
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 1.826566
$arr[12,13] = 5.479698
$arr[12,14] = 0.02795372904983374
$arr[12,15] = 0.02795372904983374
$arr[12,16] = 96.49861607748601
$arr[12,17] = 868.4875446973741
$arr[12,18] = 0.0009579247756614584
$arr[12,19] = 0.0009579247756614584
$arr[13,0] = "MuSCs"
$arr[13,1] = "Fn1"
$arr[13,2] = "Sdc2"
$arr[13,3] = "FAPs"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 52.83062100000001
$arr[13,7] = 158.491863
$arr[13,8] = 0.0342682285413064
$arr[13,9] = 0.03426822854130639
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 44.29005966666667
$arr[13,13] = 132.870179
$arr[13,14] = 0.6778141756295529
$arr[13,15] = 0.6778141756295529
$arr[13,16] = 2339.871356317054
$arr[13,17] = 21058.84220685348
$arr[13,18] = 0.02322749107901071
$arr[13,19] = 0.02322749107901071
$arr[14,0] = "MuSCs"
$arr[14,1] = "Fn1"
$arr[14,2] = "Sdc2"
$arr[14,3] = "MuSCs"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 52.83062100000001
$arr[14,7] = 158.491863
$arr[14,8] = 0.0342682285413064
$arr[14,9] = 0.03426822854130639
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 19.10886933333333
$arr[14,13] = 57.326608
$arr[14,14] = 0.2924417490485847
$arr[14,15] = 0.2924417490485847
$arr[14,16] = 1009.533433487856
$arr[14,17] = 9085.800901390705
$arr[14,18] = 0.01002146069141627
$arr[14,19] = 0.01002146069141627
$arr[15,0] = "MuSCs"
$arr[15,1] = "Fn1"
$arr[15,2] = "Sdc2"
$arr[15,3] = "Resolving-Mac"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 52.83062100000001
$arr[15,7] = 158.491863
$arr[15,8] = 0.0342682285413064
$arr[15,9] = 0.03426822854130639
$arr[15,10] = 2
$arr[15,11] = 0.6666666666666666
$arr[15,12] = 0.1169856666666667
$arr[15,13] = 0.350957
$arr[15,14] = 0.001790346272028586
$arr[15,15] = 0.001790346272028586
$arr[15,16] = 6.180425418099
$arr[15,17] = 55.623828762891
$arr[15,18] = 0.00006135199521795151
$arr[15,19] = 0.00006135199521795151
$arr[16,0] = "Resolving-Mac"
$arr[16,1] = "Fn1"
$arr[16,2] = "Sdc2"
$arr[16,3] = "ECs"
$arr[16,4] = 3
$arr[16,5] = 1
$arr[16,6] = 16.16161433333333
$arr[16,7] = 48.484843
$arr[16,8] = 0.01048312291409786
$arr[16,9] = 0.01048312291409786
$arr[16,10] = 3
$arr[16,11] = 1
$arr[16,12] = 1.826566
$arr[16,13] = 5.479698
$arr[16,14] = 0.02795372904983374
$arr[16,15] = 0.02795372904983374
$arr[16,16] = 29.52025524637933
$arr[16,17] = 265.682297217414
$arr[16,18] = 0.0002930423775367952
$arr[16,19] = 0.0002930423775367952
$arr[17,0] = "Resolving-Mac"
$arr[17,1] = "Fn1"
$arr[17,2] = "Sdc2"
$arr[17,3] = "FAPs"
$arr[17,4] = 3
$arr[17,5] = 1
$arr[17,6] = 16.16161433333333
$arr[17,7] = 48.484843
$arr[17,8] = 0.01048312291409786
$arr[17,9] = 0.01048312291409786
$arr[17,10] = 3
$arr[17,11] = 1
$arr[17,12] = 44.29005966666667
$arr[17,13] = 132.870179
$arr[17,14] = 0.6778141756295529
$arr[17,15] = 0.6778141756295529
$arr[17,16] = 715.7988631329886
$arr[17,17] = 6442.189768196897
$arr[17,18] = 0.007105609316042521
$arr[17,19] = 0.007105609316042519
$arr[18,0] = "Resolving-Mac"
$arr[18,1] = "Fn1"
$arr[18,2] = "Sdc2"
$arr[18,3] = "MuSCs"
$arr[18,4] = 3
$arr[18,5] = 1
$arr[18,6] = 16.16161433333333
$arr[18,7] = 48.484843
$arr[18,8] = 0.01048312291409786
$arr[18,9] = 0.01048312291409786
$arr[18,10] = 3
$arr[18,11] = 1
$arr[18,12] = 19.10886933333333
$arr[18,13] = 57.326608
$arr[18,14] = 0.2924417490485847
$arr[18,15] = 0.2924417490485847
$arr[18,16] = 308.8301765113938
$arr[18,17] = 2779.471588602544
$arr[18,18] = 0.003065702800490076
$arr[18,19] = 0.003065702800490075
$arr[19,0] = "Resolving-Mac"
$arr[19,1] = "Fn1"
$arr[19,2] = "Sdc2"
$arr[19,3] = "Resolving-Mac"
$arr[19,4] = 3
$arr[19,5] = 1
$arr[19,6] = 16.16161433333333
$arr[19,7] = 48.484843
$arr[19,8] = 0.01048312291409786
$arr[19,9] = 0.01048312291409786
$arr[19,10] = 2
$arr[19,11] = 0.6666666666666666
$arr[19,12] = 0.1169856666666667
$arr[19,13] = 0.350957
$arr[19,14] = 0.001790346272028586
$arr[19,15] = 0.001790346272028586
$arr[19,16] = 1.890677227194555
$arr[19,17] = 17.016095044751
$arr[19,18] = 0.00001876842002847256
$arr[19,19] = 0.00001876842002847256
$ws.Range("A2:T21").Value2 = $arr
